# Tower_logs_metadata.xlsx — "Updated metadata injection spreadsheet"
#
# The data source files referenced by this metadata workbook were renamed
# on disk (plain names -> names carrying their literal ".csv" extension),
# so the tab names that mirror those source file names are updated to
# match: callrecords_tower1/2/3  ->  callrecords_tower1/2/3.csv

$wb = $excel.ActiveWorkbook

$wb.Worksheets.Item("callrecords_tower1").Name = "callrecords_tower1.csv"
$wb.Worksheets.Item("callrecords_tower2").Name = "callrecords_tower2.csv"
$wb.Worksheets.Item("callrecords_tower3").Name = "callrecords_tower3.csv"
